$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right by one.
$ws.Range("A1").EntireColumn.Insert()

# New column A header (row 2) - "Match ID" label, bold/no-border style matching the
# other data cells in that column.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# New column A data values (rows 4-19) - all matches are Match ID 7.
$ws.Range("A4:A19").Value = 7
$ws.Range("A4:A19").Font.Bold = $true

# Row 3 is an empty spacer row (hidden) - just needs the bold style on the new cell,
# no value. Toggle visibility around the write so the runtime doesn't recompute a
# custom row height for a previously cell-less hidden row.
$ws.Rows(3).Hidden = $false
$ws.Range("A3").Font.Bold = $true
$ws.Rows(3).Hidden = $true

# Row 20 is the hidden "16 Players" summary row - gets the Match ID value but keeps
# the default (unstyled) look. Same visibility toggle trick as row 3.
$ws.Rows(20).Hidden = $false
$ws.Range("A20").Value = 7
$ws.Rows(20).Hidden = $true

# Update the sheet selection to reflect the newly inserted Match ID column.
$ws.Range("A2:A19").Select()

Write-Output "done"
